# Remove the standalone italic "Nehemia" paragraph that sits between the
# "NEH" Heading2 paragraph and the following (mostly blank) paragraph.
# (There is a later, separate, non-italic "Nehemia" Heading2 paragraph
# further down in the document that must NOT be touched.)

$d = $word.ActiveDocument

$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.Trim()
    if ($text -eq "Nehemia" -and $p.Range.Italic -eq -1) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
